# Daily attendance processing - 2025-11-25 16:32:24
# Normalize "Recorded By" (column G) values so that the System
# recorder is listed last instead of first, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
